$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Illes Balears" (A25) and "Illes Balears*" (A26)
$ws.Range("A25").Value = "Illes Balears*"
$ws.Range("A26").Value = "Illes Balears"

# Swap "Melilla" (A52) and "Huelva" (A53)
$ws.Range("A52").Value = "Huelva"
$ws.Range("A53").Value = "Melilla"

# Update the "Datos actualizados..." timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:30"
